$wb = $excel.ActiveWorkbook

# --- Sheet: Trends Status ---
$ws = $wb.Worksheets.Item("Trends Status")
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 1.9
$ws.Range("E2").Value = 6

$ws.Range("C3").Value = 23
$ws.Range("D3").Value = 12.6

$ws.Range("B4").Value = 67
$ws.Range("C4").Value = 107
$ws.Range("D4").Value = 65
$ws.Range("E4").Value = 71.3

$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 14.6
$ws.Range("E5").Value = 4

$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5.8
$ws.Range("E6").Value = 3.3

$ws.Range("B7").Value = 203
$ws.Range("C7").Value = 294

# --- Sheet: Species qualification ---
$ws = $wb.Worksheets.Item("Species qualification")
$ws.Range("C3").Value = 103
$ws.Range("C4").Value = 150

# --- Sheet: Interannual update - High Pri ---
$ws = $wb.Worksheets.Item("Interannual update - High Pri")
$ws.Range("B2").Value = 16
$ws.Range("C2").Value = 15.5
$ws.Range("D2").Value = 16
$ws.Range("E2").Value = 53.3

$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 8.699999999999999
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 10

$ws.Range("B4").Value = 78
$ws.Range("C4").Value = 75.7
$ws.Range("D4").Value = 11
$ws.Range("E4").Value = 36.7
